$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Agt"
$ws.Range("C2").Value = "Agtr2"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 3.614029333333333
$ws.Range("H2").Value = 10.842088
$ws.Range("I2").Value = 0.9241192353022275
$ws.Range("J2").Value = 0.9439348023377215
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.743532666666667
$ws.Range("N2").Value = 11.230598
$ws.Range("O2").Value = 0.9802973346235675
$ws.Range("P2").Value = 0.9802973346235675
$ws.Range("Q2").Value = 13.52923686762489
$ws.Range("R2").Value = 121.763131808624
$ws.Range("S2").Value = 0.905911623241143
$ws.Range("T2").Value = 0.9253367707900925

$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Agt"
$ws.Range("C3").Value = "Agtr2"
$ws.Range("D3").Value = "M2"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 3.614029333333333
$ws.Range("H3").Value = 10.842088
$ws.Range("I3").Value = 0.9241192353022275
$ws.Range("J3").Value = 0.9439348023377215
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.07524
$ws.Range("N3").Value = 0.22572
$ws.Range("O3").Value = 0.01970266537643246
$ws.Range("P3").Value = 0.01970266537643246
$ws.Range("Q3").Value = 0.27191956704
$ws.Range("R3").Value = 2.44727610336
$ws.Range("S3").Value = 0.01820761206108444
$ws.Range("T3").Value = 0.01859803154762904

$ws.Range("A4").Value = "M1"
$ws.Range("B4").Value = "Agt"
$ws.Range("C4").Value = "Agtr2"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.05046166666666666
$ws.Range("H4").Value = 0.151385
$ws.Range("I4").Value = 0.01290321480846011
$ws.Range("J4").Value = 0.01317989395141378
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.743532666666667
$ws.Range("N4").Value = 11.230598
$ws.Range("O4").Value = 0.9802973346235675
$ws.Range("P4").Value = 0.9802973346235675
$ws.Range("Q4").Value = 0.1889048975811111
$ws.Range("R4").Value = 1.70014407823
$ws.Range("S4").Value = 0.01264898708480879
$ws.Range("T4").Value = 0.01292021491119221

$ws.Range("A5").Value = "M1"
$ws.Range("B5").Value = "Agt"
$ws.Range("C5").Value = "Agtr2"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.05046166666666666
$ws.Range("H5").Value = 0.151385
$ws.Range("I5").Value = 0.01290321480846011
$ws.Range("J5").Value = 0.01317989395141378
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.07524
$ws.Range("N5").Value = 0.22572
$ws.Range("O5").Value = 0.01970266537643246
$ws.Range("P5").Value = 0.01970266537643246
$ws.Range("Q5").Value = 0.0037967358
$ws.Range("R5").Value = 0.0341706222
$ws.Range("S5").Value = 0.0002542277236513177
$ws.Range("T5").Value = 0.000259679040221572

$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Agt"
$ws.Range("C6").Value = "Agtr2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.5
$ws.Range("G6").Value = 0.2462915
$ws.Range("H6").Value = 0.492583
$ws.Range("I6").Value = 0.06297754988931244
$ws.Range("J6").Value = 0.04288530371086472
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.743532666666667
$ws.Range("N6").Value = 11.230598
$ws.Range("O6").Value = 0.9802973346235675
$ws.Range("P6").Value = 0.9802973346235675
$ws.Range("Q6").Value = 0.9220002757723333
$ws.Range("R6").Value = 5.532001654634001
$ws.Range("S6").Value = 0.06173672429761573
$ws.Range("T6").Value = 0.04204034892228287

$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Agt"
$ws.Range("C7").Value = "Agtr2"
$ws.Range("D7").Value = "M2"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.5
$ws.Range("G7").Value = 0.2462915
$ws.Range("H7").Value = 0.492583
$ws.Range("I7").Value = 0.06297754988931244
$ws.Range("J7").Value = 0.04288530371086472
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.07524
$ws.Range("N7").Value = 0.22572
$ws.Range("O7").Value = 0.01970266537643246
$ws.Range("P7").Value = 0.01970266537643246
$ws.Range("Q7").Value = 0.01853097246
$ws.Range("R7").Value = 0.11118583476
$ws.Range("S7").Value = 0.001240825591696704
$ws.Range("T7").Value = 0.0008449547885818448
